$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$metadata = $wb.Worksheets.Item("Metadata")

# Version: 0.9.0 -> 1.0.0
$metadata.Range("B3").Value = "1.0.0"

# Date: 2025-02-15T12:00:41+01:00 -> 2025-06-05T14:31:57+02:00
$metadata.Range("B8").Value = "2025-06-05T14:31:57+02:00"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Clear the "Condition(s)" entries that pointed to the ele-1 constraint
$elements.Range("AI4").Value = ""
$elements.Range("AI6").Value = ""
